# Updates the cryptos list (Price / Volume(1h) columns, plus the
# Avalanche/ShibaInu row swap) to match the latest scrape.
# NOTE: Price-column values that look like plain numbers (e.g. "188.64")
# are written with a leading apostrophe so Excel keeps them as literal
# text (matching the source data's original string type) instead of
# auto-converting them to floating-point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '74.883.92'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '2.849.34'
$ws.Range('E3').Value = '  +10.21%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '''601.25'
$ws.Range('E5').Value = '  +3.87%  '
$ws.Range('D6').Value = '''188.64'
$ws.Range('E6').Value = '  +2.09%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '''0.556'
$ws.Range('E8').Value = '  +4.44%  '
$ws.Range('D9').Value = '''0.194'
$ws.Range('E9').Value = '  -4.49%  '
$ws.Range('D10').Value = '2.846.17'
$ws.Range('E10').Value = '  +10.00%  '
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').Value = '''0.372'
$ws.Range('E12').Value = '  +3.95%  '
$ws.Range('D13').Value = '''4.90'
$ws.Range('E13').Value = '  +2.88%  '
$ws.Range('D14').Value = '3.370.36'
$ws.Range('E14').Value = '  +10.62%  '
$ws.Range('D15').Value = '74.840.49'
$ws.Range('E15').Value = '  +0.93%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '''0.0000188'
$ws.Range('E16').Value = '  -1.12%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').Value = '''27.25'
$ws.Range('E17').Value = '  +4.51%  '
$ws.Range('D18').Value = '2.846.69'
$ws.Range('E18').Value = '  +10.59%  '
$ws.Range('E19').Value = '  +6.09%  '
$ws.Range('D20').Value = '''12.48'
$ws.Range('E20').Value = '  +6.87%  '
$ws.Range('D21').Value = '''376.13'
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('D22').Value = '''2.27'
$ws.Range('E22').Value = '  -1.87%  '
$ws.Range('D23').Value = '''4.14'
$ws.Range('E23').Value = '  +2.61%  '
$ws.Range('D24').Value = '''6.22'
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').Value = '''0.999'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').Value = '''70.86'
$ws.Range('E26').Value = '  +1.89%  '
$ws.Range('D27').Value = '''4.22'
$ws.Range('E27').Value = '  +2.26%  '
$ws.Range('E28').Value = '  +9.82%  '
$ws.Range('D29').Value = '''9.61'
$ws.Range('E29').Value = '  +5.39%  '
$ws.Range('E30').Value = '  +11.39%  '
$ws.Range('E31').Value = '  +0.42%  '
$ws.Range('D32').Value = '''530.04'
$ws.Range('E32').Value = '  +6.15%  '
$ws.Range('D33').Value = '''1.41'
$ws.Range('E33').Value = '  +4.96%  '
$ws.Range('D34').Value = '''7.94'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('E35').Value = '  +6.54%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = '''0.121'
$ws.Range('E37').Value = '  +2.15%  '
$ws.Range('D38').Value = '''20.27'
$ws.Range('E38').Value = '  +5.99%  '
$ws.Range('D39').Value = '''162.08'
$ws.Range('E39').Value = '  +1.42%  '
$ws.Range('D40').Value = '''19.30'
$ws.Range('E40').Value = '  -0.38%  '
$ws.Range('D41').Value = '''186.33'
$ws.Range('E41').Value = '  +25.49%  '
$ws.Range('D43').Value = '''5.09'
$ws.Range('E43').Value = '  +3.75%  '
$ws.Range('E44').Value = '  +6.69%  '
$ws.Range('E45').Value = '  +1.75%  '
$ws.Range('E46').Value = '  +7.98%  '
$ws.Range('D47').Value = '''39.57'
$ws.Range('E47').Value = '  +1.84%  '
$ws.Range('D48').Value = '''2.37'
$ws.Range('E48').Value = '  -2.52%  '
$ws.Range('D49').Value = '''0.0857'
$ws.Range('E49').Value = '  +3.48%  '
$ws.Range('D50').Value = '''0.576'
$ws.Range('E50').Value = '  +11.08%  '
$ws.Range('E51').Value = '  +4.38%  '
